$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force the engine to coalesce adjacent, identically-formatted runs
# that span a given [start,end) character range by performing a harmless
# insert+delete text mutation at $pos. This mirrors what Word itself does
# when the user edits text that happens to touch several runs: it collapses
# any runs that no longer need to stay separate.
# ---------------------------------------------------------------------------
function Merge-RunsAt($doc, $pos) {
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter("X")
    $r2 = $doc.Range($pos, $pos + 1)
    $r2.Delete()
}

# ===========================================================================
# Change 1: paragraph "Testy boli vykonané pre formulárové časti systému..."
# Several runs ("Test", "y ", "bol", "i", " vykonan", "é", ...) get merged
# into a single run. No actual text changes, only run structure.
# ===========================================================================
$p1 = $d.Paragraphs.Item(7)
$p1Start = $p1.Range.Start
Merge-RunsAt $d $p1Start

# ===========================================================================
# Change 2: paragraph "Testovala sa správnosť výpočtov. ..."
# The first run (big sentence) stays as-is; the trailing three runs
# (" ", "Systém obsahuje...opraví", ". Tento algoritmus...testovaný.")
# get merged into a single run.
# ===========================================================================
$p2 = $d.Paragraphs.Item(20)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1
$p2Text = $d.Range($p2Start, $p2End).Text
$boundaryOffset = $p2Text.IndexOf("hodnotu.") + "hodnotu.".Length
$boundaryPos = $p2Start + $boundaryOffset

# Protect the boundary between the first (unchanged) run and the group of
# three runs that must merge, so the merge doesn't swallow the first run.
$protectRange = $d.Range($boundaryPos, $boundaryPos)
$d.Bookmarks.Add("TMP_PROTECT_2", $protectRange) | Out-Null

Merge-RunsAt $d ($boundaryPos + 1)

$d.Bookmarks.Item("TMP_PROTECT_2").Delete()

# ===========================================================================
# Change 3: paragraph "Komponent podľa predpísaných postupov ..."
# The run "vylepšenia sú presnejšie definované v dokumente „Vlastné
# zhodnotenie diela a návrhy na vylepšienie“." is split in three, a
# stray "i" is removed (vylepšienie -> vylepšenie), and the _GoBack
# bookmark ends up between the 2nd and 3rd pieces.
# ===========================================================================
$p3 = $d.Paragraphs.Item(32)
$p3Start = $p3.Range.Start
$p3End = $p3.Range.End - 1
$p3Text = $d.Range($p3Start, $p3End).Text

# Boundary between " Tieto " and "vylepšenia sú ..." - must be preserved.
$tietoBoundaryOffset = $p3Text.IndexOf("vylepšenia sú")
$tietoBoundaryPos = $p3Start + $tietoBoundaryOffset

# Split point 1: right after "...Vlastné zhodn" (before "otenie...").
$split1Offset = $p3Text.IndexOf("zhodnotenie") + "zhodn".Length
$split1Pos = $p3Start + $split1Offset

# Split point 2: right after the final "...na vylepš" (before "ienie...").
$split2Offset = $p3Text.LastIndexOf("vylepš") + "vylepš".Length
$split2Pos = $p3Start + $split2Offset

# Step A: protect the " Tieto " / "vylepšenia" boundary while we edit text
# further along in the same run (otherwise the engine would coalesce the
# two runs together since they share identical formatting).
$d.Bookmarks.Add("TMP_PROTECT_3", $d.Range($tietoBoundaryPos, $tietoBoundaryPos)) | Out-Null

# Step B: delete the stray "i" so "vylepšienie" becomes "vylepšenie".
$iCharRange = $d.Range($split2Pos, $split2Pos + 1)
$iCharRange.Delete()

# Step C: boundary no longer needs protecting.
$d.Bookmarks.Item("TMP_PROTECT_3").Delete()

# Step D: split the run at split1Pos (between "zhodn" and "otenie").
$d.Bookmarks.Add("TMP_SPLIT_3", $d.Range($split1Pos, $split1Pos)) | Out-Null
$d.Bookmarks.Item("TMP_SPLIT_3").Delete()

# Step E: move the _GoBack bookmark to split2Pos (between "vylepš" and "enie").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($split2Pos, $split2Pos)) | Out-Null
